# Update the TestResult sheet: populate the Status column (B) for each
# test script row with a PASS/FAIL result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResult")

$ws.Range("B2").Value = "PASS"
$ws.Range("B3").Value = "FAIL"
$ws.Range("B4").Value = "PASS"
